$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.852.85"
$ws.Range("E2").Value = "  +5.80%  "
$ws.Range("D3").Value = "3.537.21"
$ws.Range("E3").Value = "  +8.91%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'566.93"
$ws.Range("E5").Value = "  +7.03%  "
$ws.Range("D6").Value = "'188.74"
$ws.Range("E6").Value = "  +9.96%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'0.619"
$ws.Range("E7").Value = "  +3.75%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.525.84"
$ws.Range("E8").Value = "  +8.68%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "'0.634"
$ws.Range("E10").Value = "  +4.63%  "
$ws.Range("D11").Value = "'0.151"
$ws.Range("E11").Value = "  +12.80%  "
$ws.Range("D12").Value = "'54.84"
$ws.Range("E12").Value = "  +3.09%  "
$ws.Range("D13").Value = "'0.0000270"
$ws.Range("E13").Value = "  +6.31%  "
$ws.Range("D14").Value = "'9.41"
$ws.Range("E14").Value = "  +2.95%  "
$ws.Range("D15").Value = "4.106.22"
$ws.Range("E15").Value = "  +9.23%  "
$ws.Range("D16").Value = "3.541.75"
$ws.Range("E16").Value = "  +9.81%  "
$ws.Range("E17").Value = "  +4.41%  "
$ws.Range("D18").Value = "66.921.24"
$ws.Range("E18").Value = "  +5.98%  "
$ws.Range("D19").Value = "'18.27"
$ws.Range("E19").Value = "  +5.87%  "
$ws.Range("D20").Value = "'12.01"
$ws.Range("E20").Value = "  +8.08%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("E21").Value = "  +2.76%  "
$ws.Range("D22").Value = "'428.88"
$ws.Range("E22").Value = "  +16.87%  "
$ws.Range("D23").Value = "'4.21"
$ws.Range("E23").Value = "  +11.97%  "
$ws.Range("D24").Value = "'85.14"
$ws.Range("E24").Value = "  +5.04%  "
$ws.Range("D25").Value = "'4.13"
$ws.Range("E25").Value = "  +3.55%  "
$ws.Range("D26").Value = "'11.14"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("E27").Value = "  +9.77%  "
$ws.Range("D28").Value = "'12.27"
$ws.Range("E28").Value = "  +8.61%  "
$ws.Range("D29").Value = "'9.25"
$ws.Range("E29").Value = "  +12.69%  "
$ws.Range("D30").Value = "'30.35"
$ws.Range("E30").Value = "  +6.39%  "
$ws.Range("D31").Value = "'642.42"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("D32").Value = "'6.62"
$ws.Range("E32").Value = "  +2.40%  "
$ws.Range("D33").Value = "'11.75"
$ws.Range("E33").Value = "  +4.82%  "
$ws.Range("D34").Value = "'0.112"
$ws.Range("E34").Value = "  +6.04%  "
$ws.Range("D35").Value = "'59.89"
$ws.Range("E35").Value = "  +5.27%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").Value = "'38.66"
$ws.Range("E36").Value = "  +5.32%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.149"
$ws.Range("E37").Value = "  +21.13%  "
$ws.Range("D38").Value = "0.0₃0815"
$ws.Range("E38").Value = "  +13.50%  "
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D40").Value = "'0.393"
$ws.Range("E40").Value = "  +3.99%  "
$ws.Range("D41").Value = "'3.36"
$ws.Range("E41").Value = "  +14.51%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").Value = "3.053.77"
$ws.Range("E43").Value = "  +6.04%  "
$ws.Range("D44").Value = "'2.67"
$ws.Range("E44").Value = "  +4.62%  "
$ws.Range("D45").Value = "'2.90"
$ws.Range("E45").Value = "  +11.90%  "
$ws.Range("D46").Value = "'3.36"
$ws.Range("E46").Value = "  +9.37%  "
$ws.Range("E47").Value = "  +6.81%  "
$ws.Range("D48").Value = "'2.77"
$ws.Range("E48").Value = "  +3.29%  "
$ws.Range("E49").Value = "  +5.49%  "
$ws.Range("D50").Value = "'143.43"
$ws.Range("E50").Value = "  +6.76%  "
$ws.Range("E51").Value = "  +10.39%  "
